# Update SCH to A4 page size
# The BOM worksheet gets a new header row inserted at the top, and the
# previously-selected cell moves from A3 to A6 (net effect of the header
# row shifting all existing data down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above the current row 1; this shifts all
# existing BOM rows (previously rows 1-15) down to rows 2-16 while
# preserving their shared-string references and cell formatting.
$ws.Rows.Item(1).Insert()

# Populate the new header row. The assignment order below matches the
# order the new strings were appended to the shared-string table in the
# saved workbook (Designator, Value, Description, Manufacturer PN,
# Quantity, Datasheet).
$ws.Range("A1").Value = "Designator"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Description"
$ws.Range("E1").Value = "Manufacturer PN"
$ws.Range("F1").Value = "Quantity"
$ws.Range("D1").Value = "Datasheet"

# Restore the selected cell to where it ends up after the shift (the
# previously-selected A3 now corresponds to A6).
$ws.Range("A6").Select()
